$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.330.78"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").Value = "1.857.38"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("E4").Value = "  -0.70%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "314.29"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.25%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4613"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.17%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3705"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.64%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07324"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.01%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.8825"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +2.45%  "
$ws.Range("E11").Value = "  +0.72%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "19.83"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "1.888.02"
$ws.Range("E13").Value = "  +2.13%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.385"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.89%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "6.559"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.34%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "91.92"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.15%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.61%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008843"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("E19").Value = "  -0.63%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.83"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("D21").Value = "27.350.52"
$ws.Range("E21").Value = "  +0.71%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.120"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.55%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.52"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").Value = "2.125.12"
$ws.Range("E24").Value = "  +1.96%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.888"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.47%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "152.42"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("E27").Value = "  +1.02%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.083"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("E29").Value = "  +0.32%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "116.01"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.51%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.08864"
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.7613"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +5.38%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.022"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.18%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.172"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.59%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "4.487"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +1.11%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.626"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +7.25%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.01961"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.00%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.074"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.59%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.984"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.92%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.05208"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.56%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "7.045"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.32%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.5163"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("E43").Value = "  +0.72%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "8.346"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +2.02%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.4832"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.59%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "10.34"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.44%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.66%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "103.31"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.04%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.06231"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.57%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "65.73"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.15%  "
